# Applies the November 2016 EIA Table A.8.A update:
#  - subtitle month text October -> November
#  - revised relative-standard-error figures for rows 4-65 (cols B-F)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table_A_8_A")

# Title / subtitle (A2): "October 2016" -> "November 2016"
$ws.Range("A2").Value = "by End-Use Sector, Census Division, and State, November 2016"

# Revised figures
# Row 4: B4, C4, F4
$ws.Range("B4").Value = 0.28000000000000003
$ws.Range("C4").Value = 0.24
$ws.Range("F4").Value = 0.36

# Row 5: B5, C5, F5
$ws.Range("B5").Value = 0.28000000000000003
$ws.Range("C5").Value = 0.43
$ws.Range("F5").Value = 0.49

# Row 6: B6, C6
$ws.Range("B6").Value = 0.41
$ws.Range("C6").Value = 0.45

# Row 7: C7
$ws.Range("C7").Value = 0.41

# Row 8: B8, C8
$ws.Range("B8").Value = 0.4
$ws.Range("C8").Value = 0.41

# Row 10: C10
$ws.Range("C10").Value = 2

# Row 11: B11, C11, D11
$ws.Range("B11").Value = 0.14000000000000001
$ws.Range("C11").Value = 0.09
$ws.Range("D11").Value = 0.47

# Row 12: B12, C12, F12
$ws.Range("B12").Value = 0.26
$ws.Range("C12").Value = 0.2
$ws.Range("F12").Value = 0.2

# Row 13: B13, C13, F13
$ws.Range("B13").Value = 0.17
$ws.Range("C13").Value = 0.11
$ws.Range("F13").Value = 0.24

# Row 14: B14, C14, F14
$ws.Range("B14").Value = 0.28999999999999998
$ws.Range("C14").Value = 0.21
$ws.Range("F14").Value = 0.19

# Row 15: B15, C15, F15
$ws.Range("B15").Value = 0.3
$ws.Range("C15").Value = 0.3
$ws.Range("F15").Value = 0.23

# Row 16: C16, D16, F16
$ws.Range("C16").Value = 0.44
$ws.Range("D16").Value = 1
$ws.Range("F16").Value = 0.37

# Row 18: F18
$ws.Range("F18").Value = 0.41

# Row 19: C19
$ws.Range("C19").Value = 0.41

# Row 21: F21
$ws.Range("F21").Value = 0.45

# Row 23: B23, C23
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 2

# Row 25: D25
$ws.Range("D25").Value = 2

# Row 29: B29, C29, F29
$ws.Range("B29").Value = 0.31
$ws.Range("C29").Value = 0.39
$ws.Range("F29").Value = 0.25

# Row 30: D30
$ws.Range("D30").Value = 6

# Row 32: B32, F32
$ws.Range("B32").Value = 0.39
$ws.Range("F32").Value = 0.38

# Row 34: B34, C34, D34
$ws.Range("B34").Value = 0.32
$ws.Range("C34").Value = 0.27
$ws.Range("D34").Value = 1

# Row 35: D35
$ws.Range("D35").Value = 2

# Row 36: D36
$ws.Range("D36").Value = 2

# Row 37: F37
$ws.Range("F37").Value = 0.43

# Row 38: B38, C38, D38, F38
$ws.Range("B38").Value = 0.26
$ws.Range("C38").Value = 0.35
$ws.Range("D38").Value = 0.12
$ws.Range("F38").Value = 0.16

# Row 39: F39
$ws.Range("F39").Value = 0.44

# Row 41: D41
$ws.Range("D41").Value = 1

# Row 42: B42, D42
$ws.Range("B42").Value = 1
$ws.Range("D42").Value = 3

# Row 44: B44, F44
$ws.Range("B44").Value = 0.49
$ws.Range("F44").Value = 0.37

# Row 47: C47
$ws.Range("C47").Value = 2

# Row 48: F48
$ws.Range("F48").Value = 0.41

# Row 49: F49
$ws.Range("F49").Value = 0.37

# Row 50: B50, F50
$ws.Range("B50").Value = 1
$ws.Range("F50").Value = 1

# Row 53: B53, D53
$ws.Range("B53").Value = 1
$ws.Range("D53").Value = 9

# Row 54: B54, D54
$ws.Range("B54").Value = 0.41
$ws.Range("D54").Value = 2

# Row 55: D55
$ws.Range("D55").Value = 6

# Row 56: C56
$ws.Range("C56").Value = 1

# Row 58: B58, F58
$ws.Range("B58").Value = 0.24
$ws.Range("F58").Value = 0.36

# Row 59: B59, C59, F59
$ws.Range("B59").Value = 0.2
$ws.Range("C59").Value = 0.37
$ws.Range("F59").Value = 0.27

# Row 63: B63, D63
$ws.Range("B63").Value = 1
$ws.Range("D63").Value = 5

# Row 65: B65, C65, D65
$ws.Range("B65").Value = 0.17
$ws.Range("C65").Value = 0.19
$ws.Range("D65").Value = 0.42
